# Apply the "room schedule" header block (room name / day / AP info) that
# the Report app now generates: B-203B room header, Monday column header,
# and the access-point info row (A1 / ADSL / SEA / APK) inserted into the
# existing 8:00-9:00 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors expressed as OLE COLORREF (0xBBGGRR) so Interior.Color renders the
# intended RGB on save: red=FF0000, yellow=FFFF00, blue=77ABFF.
$RED    = 255        # 0x0000FF -> RGB FF0000
$YELLOW = 65535       # 0x00FFFF -> RGB FFFF00
$BLUE   = 16755575    # 0xFFAB77 -> RGB 77ABFF

$CENTER = -4108       # xlCenter

# --- Row 1: "Monday" header spanning B1:C1 (bold, red fill, boxed, centered)
$ws.Range("B1:C1").Merge()
$ws.Range("B1").Value = "Monday"
$hdr1 = $ws.Range("B1:C1")
$hdr1.Font.Bold = $true
$hdr1.Interior.Color = $RED
$hdr1.Borders.LineStyle = 1
$hdr1.HorizontalAlignment = $CENTER
$hdr1.VerticalAlignment = $CENTER

# --- Row 2: "B-203B" room header spanning B2:C2 (bold, yellow fill, boxed, centered)
$ws.Range("B2:C2").Merge()
$ws.Range("B2").Value = "B-203B"
$hdr2 = $ws.Range("B2:C2")
$hdr2.Font.Bold = $true
$hdr2.Interior.Color = $YELLOW
$hdr2.Borders.LineStyle = 1
$hdr2.HorizontalAlignment = $CENTER
$hdr2.VerticalAlignment = $CENTER

# --- Row 7/8: access point info (A1/ADSL on row7, SEA/APK on row8), blue fill
$ws.Range("B7").Value = "A1"
$ws.Range("C7").Value = "ADSL"
$ws.Range("B8").Value = "SEA"
$ws.Range("C8").Value = "APK"
$ap = $ws.Range("B7:C8")
$ap.Interior.Color = $BLUE
$ap.Borders.LineStyle = 1
$ap.HorizontalAlignment = $CENTER
$ap.VerticalAlignment = $CENTER
